# Apply updated cryptocurrency price/volume figures to sheet1 (cells D2:E51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "27.590.65"
$ws.Range('E2').Value = "  +0.27%  "

$ws.Range('D3').Value = "1.847.78"
$ws.Range('E3').Value = "  +0.25%  "

$ws.Range('D4').Value = "'1.026"
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = "  -0.47%  "

$ws.Range('D5').Value = "'321.14"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = "  +0.44%  "

$ws.Range('E6').Value = "  -0.38%  "

$ws.Range('D7').Value = "'0.4368"
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = "  -0.14%  "

$ws.Range('D8').Value = "'0.3786"
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = "  +1.29%  "

$ws.Range('D9').Value = "'0.07372"
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = "  -0.33%  "

$ws.Range('D10').Value = "'0.8802"
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = "  +0.32%  "

$ws.Range('D11').Value = "'21.47"
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = "  +0.25%  "

$ws.Range('D12').Value = "1.862.12"
$ws.Range('E12').Value = "  +0.78%  "

$ws.Range('D13').Value = "'5.491"
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = "  +0.00%  "

$ws.Range('D14').Value = "'6.707"
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = "  +0.29%  "

$ws.Range('D15').Value = "'0.07121"
$ws.Range('D15').Style = "Normal"

$ws.Range('D16').Value = "'85.05"
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = "  +2.81%  "

$ws.Range('D17').Value = "'1.032"
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = "  -0.16%  "

$ws.Range('D18').Value = "'0.000009049"
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = "  +0.05%  "

$ws.Range('D19').Value = "'1.025"
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = "  -0.30%  "

$ws.Range('D20').Value = "'15.36"
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = "  -0.29%  "

$ws.Range('D21').Value = "27.610.30"
$ws.Range('E21').Value = "  +0.25%  "

$ws.Range('D22').Value = "'5.281"
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = "  +0.94%  "

$ws.Range('D23').Value = "'11.24"
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = "  +0.43%  "

$ws.Range('D24').Value = "2.090.66"
$ws.Range('E24').Value = "  +0.79%  "

$ws.Range('D25').Value = "'2.021"
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = "  +4.78%  "

$ws.Range('D26').Value = "'157.00"
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = "  -0.06%  "

$ws.Range('D27').Value = "'18.72"
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = "  -0.14%  "

$ws.Range('D28').Value = "'5.327"
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = "  +1.41%  "

$ws.Range('D29').Value = "'1.981"
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = "  +2.15%  "

$ws.Range('D30').Value = "'117.37"
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = "  +0.90%  "

$ws.Range('E31').Value = "  -0.97%  "

$ws.Range('D32').Value = "'0.7720"
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = "  +1.17%  "

$ws.Range('D33').Value = "'1.208"
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = "  -0.18%  "

$ws.Range('E34').Value = "  +3.71%  "

$ws.Range('D35').Value = "'4.546"
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = "  +1.32%  "

$ws.Range('D36').Value = "'1.026"
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = "  -0.53%  "

$ws.Range('E37').Value = "  -0.90%  "

$ws.Range('D38').Value = "'0.01969"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = "  -0.03%  "

$ws.Range('D39').Value = "'0.05255"
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = "  +0.03%  "

$ws.Range('D40').Value = "'2.842"
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = "  +2.03%  "

$ws.Range('D41').Value = "'0.5176"
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = "  +0.03%  "

$ws.Range('E42').Value = "  +0.70%  "

$ws.Range('D43').Value = "'6.840"
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = "  +3.31%  "

$ws.Range('D44').Value = "'8.772"
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = "  +3.06%  "

$ws.Range('D45').Value = "'109.95"
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = "  +0.75%  "

$ws.Range('D46').Value = "'10.62"
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = "  +0.66%  "

$ws.Range('E47').Value = "  +3.99%  "

$ws.Range('D48').Value = "'1.026"
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = "  -0.41%  "

$ws.Range('D49').Value = "'1.700"
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = "  -0.24%  "

$ws.Range('D50').Value = "'0.4688"
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = "  +1.06%  "

$ws.Range('D51').Value = "'1.898"
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = "  +0.88%  "
